$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 (old "test2@gmail.com" / "Test 2" row) is being cleared out ---
# Remember B3's value is removed entirely, A3 becomes a blank (but still
# carries its existing "Hyperlink-look" style s=1) cell with no hyperlink.
$ws.Range("B3").ClearContents()

# --- Row 2: the email address on A2 is replaced; B2 keeps "Test 1" ---
$ws.Range("A2").Value = "akuntumbalknox@gmail.com"
$ws.Range("B2").Value = "Test 1"
$ws.Range("A3").ClearContents()

# --- Hyperlinks: drop both existing links (A2 -> test1@, A3 -> test2@) and
# recreate only the one that should remain, now pointing at the new address.
# (This runtime only supports wiping the whole collection, not per-item
# deletion, so we rebuild the single surviving link from scratch.)
$ws.Range("A2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A2"), "mailto:akuntumbalknox@gmail.com")

# Re-applying the Hyperlink cell style keeps A2 on the same style index (s=1)
# it already had before -- Hyperlinks.Add() would otherwise stamp a fresh
# (duplicate) style record onto the cell.
$ws.Range("A2").Style = "Hyperlink"

# --- New trailing blank row (row 5), styled like the existing blank row 4 ---
$ws.Range("A4").Copy()
$ws.Range("A5").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Selection / view bookkeeping to match the saved workbook state ---
$ws.Range("B3").Select()
$excel.ActiveWindow.Left = 4884
$excel.ActiveWindow.Top = 3060
$excel.ActiveWindow.Width = 17280
$excel.ActiveWindow.Height = 8880
